$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Update the raw measured test data (columns E and F) ---
# Column E: "Time for 5000 iterations (s)"; Column F: "Iterations per second" (= 5000 / E)
$ws.Range("E4").Value = 188.28299999999999
$ws.Range("F4").Value = 26.555700000000002

$ws.Range("E5").Value = 332.30099999999999
$ws.Range("F5").Value = 15.0466

$ws.Range("E6").Value = 1356.74
$ws.Range("F6").Value = 3.6853199999999999

$ws.Range("E7").Value = 169.077
$ws.Range("F7").Value = 29.572299999999998

$ws.Range("E8").Value = 970.74800000000005
$ws.Range("F8").Value = 5.1506699999999999

$ws.Range("E9").Value = 956.35599999999999
$ws.Range("F9").Value = 5.22818

# --- Resize / move the embedded chart (its bottom-right anchor shifted) ---
$co = $ws.ChartObjects(1)
$co.Width = 324.84385826771654
$co.Height = 311.62503937007875

# --- Update the selection shown in the sheet view ---
$ws.Range("A3:F9").Select()
